$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 302
$endRow = 328
$startDate = 44376

# Copy the formatting (style) of the last existing row's A cell (A301)
# down into the new A cells we are about to create.
$formatSource = $ws.Range("A301")

for ($r = $startRow; $r -le $endRow; $r++) {
    $dateSerial = $startDate + ($r - $startRow)

    $formatSource.Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
